# Replace the multiplication problems in the table with new values,
# per the commit diff. Each old value is unique within the document,
# so a straightforward Find/Replace (MatchWholeWord) per pair is safe.

$d = $word.ActiveDocument

$pairs = @(
    @("33×63=2079", "87×37=3219"),
    @("53×28=1484", "27×32=864"),
    @("63×54=3402", "42×84=3528"),
    @("36×74=2664", "58×86=4988"),
    @("25×79=1975", "61×86=5246"),
    @("28×19=532",  "22×42=924"),
    @("71×81=5751", "69×54=3726"),
    @("95×55=5225", "33×17=561"),
    @("89×58=5162", "26×85=2210"),
    @("18×28=504",  "87×33=2871"),
    @("39×40=1560", "52×23=1196"),
    @("64×75=4800", "52×57=2964"),
    @("83×72=5976", "56×71=3976"),
    @("51×47=2397", "64×97=6208"),
    @("17×93=1581", "69×28=1932"),
    @("61×77=4697", "68×70=4760"),
    @("29×97=2813", "67×68=4556"),
    @("85×48=4080", "11×33=363"),
    @("47×96=4512", "98×45=4410"),
    @("48×92=4416", "77×15=1155"),
    @("52×84=4368", "14×92=1288"),
    @("90×97=8730", "60×97=5820"),
    @("70×39=2730", "34×50=1700"),
    @("54×23=1242", "32×31=992"),
    @("85×71=6035", "21×47=987")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
